$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: columns get re-labelled (and their data columns move along with the label)
$ws.Range("B1").Value = "CA_LF"
$ws.Range("C1").Value = "FFR_LF"
$ws.Range("D1").Value = "LF_CA"
$ws.Range("E1").Value = "LF_FFR"

# Row 2 ("params") - values follow their (re-ordered) labels, with refreshed precision
$ws.Range("B2").Value = 0.9168477780017975
$ws.Range("C2").Value = 1.406551190655028
$ws.Range("D2").Value = 0.1382174556653223
$ws.Range("E2").Value = 0.6148321271736868

# Row 3 ("pvalue") - values follow their (re-ordered) labels
$ws.Range("B3").Value = 0.000001486500189606943
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.0000000181484669514731
$ws.Range("E3").Value = 0
